$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 4189.75
$ws.Cells.Item(32, 9).Value = 0
$ws.Cells.Item(32, 11).Value = 0
$ws.Cells.Item(32, 13).Value = ""
$ws.Cells.Item(53, 8).Value = 2619.6667
$ws.Cells.Item(53, 9).Value = 5099.25
$ws.Cells.Item(53, 11).Value = 5099.25
$ws.Cells.Item(53, 13).Value = -4462.25
$ws.Cells.Item(125, 8).Value = 1796.5
$ws.Cells.Item(125, 10).Value = 1762
$ws.Cells.Item(125, 12).Value = 15858
$ws.Cells.Item(125, 14).Value = -20778
$ws.Cells.Item(132, 8).Value = 1253.1
$ws.Cells.Item(132, 9).Value = 1209
$ws.Cells.Item(132, 11).Value = 3627
$ws.Cells.Item(132, 13).Value = -1097
$ws.Cells.Item(135, 8).Value = 765.25
$ws.Cells.Item(135, 9).Value = 532
$ws.Cells.Item(135, 11).Value = 4788
$ws.Cells.Item(135, 13).Value = -2253
$ws.Cells.Item(138, 8).Value = 2678.6287
$ws.Cells.Item(138, 9).Value = 2974.5715
$ws.Cells.Item(138, 10).Value = 2481.3333
$ws.Cells.Item(138, 11).Value = 8923.7145
$ws.Cells.Item(138, 12).Value = 7443.999899999999
$ws.Cells.Item(138, 13).Value = -3783.7145
$ws.Cells.Item(138, 14).Value = -17723.9999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(35, 8).Value = 3418.5
$ws.Cells.Item(35, 9).Value = 3418.5
$ws.Cells.Item(35, 11).Value = 3418.5
$ws.Cells.Item(35, 13).Value = -3012.5
$ws.Cells.Item(76, 8).Value = 0
$ws.Cells.Item(76, 10).Value = 0
$ws.Cells.Item(76, 12).Value = ""
$ws.Cells.Item(76, 14).Value = 0
$ws.Cells.Item(79, 8).Value = 0
$ws.Cells.Item(79, 10).Value = 0
$ws.Cells.Item(79, 12).Value = ""
$ws.Cells.Item(79, 14).Value = 0
$ws.Cells.Item(132, 8).Value = 2336.6304
$ws.Cells.Item(132, 9).Value = 2122.4102
$ws.Cells.Item(132, 10).Value = 3530.1428
$ws.Cells.Item(132, 11).Value = 6367.230599999999
$ws.Cells.Item(132, 12).Value = 10590.4284
$ws.Cells.Item(132, 13).Value = -3837.230599999999
$ws.Cells.Item(132, 14).Value = -15650.4284

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1268.1063
$ws.Cells.Item(31, 9).Value = 833.0833
$ws.Cells.Item(31, 11).Value = 833.0833
$ws.Cells.Item(31, 13).Value = -538.0833
$ws.Cells.Item(34, 8).Value = 1268.1063
$ws.Cells.Item(34, 9).Value = 833.0833
$ws.Cells.Item(34, 11).Value = 833.0833
$ws.Cells.Item(34, 13).Value = -631.0833
$ws.Cells.Item(60, 8).Value = 4999.9
$ws.Cells.Item(60, 10).Value = 14333
$ws.Cells.Item(60, 12).Value = 14333
$ws.Cells.Item(60, 14).Value = -15355
$ws.Cells.Item(74, 8).Value = 30000
$ws.Cells.Item(74, 10).Value = 30000
$ws.Cells.Item(74, 12).Value = 30000
$ws.Cells.Item(74, 14).Value = -31748
$ws.Cells.Item(77, 8).Value = 30000
$ws.Cells.Item(77, 10).Value = 30000
$ws.Cells.Item(77, 12).Value = 90000
$ws.Cells.Item(77, 14).Value = -98736
$ws.Cells.Item(134, 8).Value = 1666.3572
$ws.Cells.Item(134, 9).Value = 875.25
$ws.Cells.Item(134, 11).Value = 2625.75
$ws.Cells.Item(134, 13).Value = -90.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 954623.8
$ws.Cells.Item(4, 9).Value = 954623.8
$ws.Cells.Item(4, 11).Value = 2863871.4
$ws.Cells.Item(4, 13).Value = -2863759.4
$ws.Cells.Item(118, 8).Value = 37039012
$ws.Cells.Item(118, 9).Value = 55556016
$ws.Cells.Item(118, 11).Value = 166668048
$ws.Cells.Item(118, 13).Value = -166666805
$ws.Cells.Item(131, 8).Value = 11129620
$ws.Cells.Item(131, 9).Value = 166667140
$ws.Cells.Item(131, 10).Value = 19796.857
$ws.Cells.Item(131, 11).Value = 500001420
$ws.Cells.Item(131, 12).Value = 59390.571
$ws.Cells.Item(131, 13).Value = -499996380
$ws.Cells.Item(131, 14).Value = -69470.571

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 430000
$ws.Cells.Item(2, 9).Value = 500000
$ws.Cells.Item(2, 10).Value = 10000
$ws.Cells.Item(2, 11).Value = 500000
$ws.Cells.Item(2, 12).Value = 10000
$ws.Cells.Item(2, 13).Value = -499888
$ws.Cells.Item(2, 14).Value = -10224
$ws.Cells.Item(3, 8).Value = 15000
$ws.Cells.Item(3, 10).Value = 15000
$ws.Cells.Item(3, 12).Value = 15000
$ws.Cells.Item(3, 14).Value = -15224
$ws.Cells.Item(4, 8).Value = 14801.8
$ws.Cells.Item(4, 9).Value = 14669.667
$ws.Cells.Item(4, 11).Value = 14669.667
$ws.Cells.Item(4, 13).Value = -14556.667
$ws.Cells.Item(10, 8).Value = 70004
$ws.Cells.Item(10, 9).Value = 0
$ws.Cells.Item(10, 11).Value = 0
$ws.Cells.Item(10, 13).Value = ""
$ws.Cells.Item(15, 8).Value = 15000
$ws.Cells.Item(15, 10).Value = 15000
$ws.Cells.Item(15, 12).Value = 15000
$ws.Cells.Item(15, 14).Value = -15340
$ws.Cells.Item(19, 8).Value = 1051.5
$ws.Cells.Item(19, 9).Value = 1051.5
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 11).Value = 1051.5
$ws.Cells.Item(19, 12).Value = 0
$ws.Cells.Item(19, 13).Value = ""
$ws.Cells.Item(19, 14).Value = -881.5
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(21, 10).Value = 0
$ws.Cells.Item(21, 12).Value = ""
$ws.Cells.Item(21, 14).Value = 0
$ws.Cells.Item(22, 8).Value = 4640
$ws.Cells.Item(22, 9).Value = 1500
$ws.Cells.Item(22, 10).Value = 6995
$ws.Cells.Item(22, 11).Value = 1500
$ws.Cells.Item(22, 12).Value = 6995
$ws.Cells.Item(22, 13).Value = -1205
$ws.Cells.Item(22, 14).Value = -7585
$ws.Cells.Item(25, 8).Value = 0
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 10).Value = 0
$ws.Cells.Item(25, 11).Value = 0
$ws.Cells.Item(25, 12).Value = ""
$ws.Cells.Item(25, 13).Value = ""
$ws.Cells.Item(25, 14).Value = 0
$ws.Cells.Item(27, 8).Value = 4640
$ws.Cells.Item(27, 9).Value = 1500
$ws.Cells.Item(27, 10).Value = 6995
$ws.Cells.Item(27, 11).Value = 1500
$ws.Cells.Item(27, 12).Value = 6995
$ws.Cells.Item(27, 13).Value = -1393
$ws.Cells.Item(27, 14).Value = -7209
$ws.Cells.Item(28, 8).Value = 14801.8
$ws.Cells.Item(28, 9).Value = 14669.667
$ws.Cells.Item(28, 11).Value = 14669.667
$ws.Cells.Item(28, 13).Value = -14437.667
$ws.Cells.Item(30, 8).Value = 3833.3333
$ws.Cells.Item(30, 9).Value = 3833.3333
$ws.Cells.Item(30, 11).Value = 3833.3333
$ws.Cells.Item(30, 13).Value = -3725.3333
$ws.Cells.Item(32, 8).Value = 4545.2
$ws.Cells.Item(32, 9).Value = 4545.2
$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 11).Value = 4545.2
$ws.Cells.Item(32, 12).Value = 0
$ws.Cells.Item(32, 13).Value = ""
$ws.Cells.Item(32, 14).Value = -4228.2
$ws.Cells.Item(35, 8).Value = 5236.2
$ws.Cells.Item(35, 9).Value = 2378.75
$ws.Cells.Item(35, 11).Value = 2378.75
$ws.Cells.Item(35, 13).Value = -2042.75
$ws.Cells.Item(37, 8).Value = 14801.8
$ws.Cells.Item(37, 9).Value = 14669.667
$ws.Cells.Item(37, 11).Value = 14669.667
$ws.Cells.Item(37, 13).Value = -14562.667
$ws.Cells.Item(43, 8).Value = 674666.7
$ws.Cells.Item(43, 10).Value = 674666.7
$ws.Cells.Item(43, 12).Value = 674666.7
$ws.Cells.Item(43, 14).Value = -675052.7
$ws.Cells.Item(45, 8).Value = 0
$ws.Cells.Item(45, 9).Value = 0
$ws.Cells.Item(45, 11).Value = 0
$ws.Cells.Item(45, 13).Value = ""
$ws.Cells.Item(46, 8).Value = 946.5454999999999
$ws.Cells.Item(46, 9).Value = 528.25
$ws.Cells.Item(46, 10).Value = 1185.5714
$ws.Cells.Item(46, 11).Value = 528.25
$ws.Cells.Item(46, 12).Value = 1185.5714
$ws.Cells.Item(46, 13).Value = -340.25
$ws.Cells.Item(46, 14).Value = -1561.5714
$ws.Cells.Item(75, 8).Value = 50000
$ws.Cells.Item(75, 10).Value = 50000
$ws.Cells.Item(75, 12).Value = 50000
$ws.Cells.Item(75, 14).Value = -51872
$ws.Cells.Item(76, 8).Value = 22000
$ws.Cells.Item(76, 10).Value = 22000
$ws.Cells.Item(76, 12).Value = 22000
$ws.Cells.Item(76, 14).Value = -22676
$ws.Cells.Item(78, 8).Value = 50000
$ws.Cells.Item(78, 10).Value = 50000
$ws.Cells.Item(78, 12).Value = 150000
$ws.Cells.Item(78, 14).Value = -159360
$ws.Cells.Item(79, 8).Value = 22000
$ws.Cells.Item(79, 10).Value = 22000
$ws.Cells.Item(79, 12).Value = 22000
$ws.Cells.Item(79, 14).Value = -24340
$ws.Cells.Item(133, 8).Value = 67883.664
$ws.Cells.Item(133, 10).Value = 67883.664
$ws.Cells.Item(133, 12).Value = 67883.664
$ws.Cells.Item(133, 14).Value = -72943.664

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(108, 8).Value = 0
$ws.Cells.Item(108, 9).Value = 0
$ws.Cells.Item(108, 10).Value = 0
$ws.Cells.Item(108, 11).Value = 0
$ws.Cells.Item(108, 12).Value = ""
$ws.Cells.Item(108, 13).Value = ""
$ws.Cells.Item(108, 14).Value = 0
$ws.Cells.Item(132, 8).Value = 3530.5
$ws.Cells.Item(132, 9).Value = 2854.4443
$ws.Cells.Item(132, 11).Value = 8563.332900000001
$ws.Cells.Item(132, 13).Value = -6033.332900000001
